$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate current row 3 (with its original values) down into a new row 4,
# preserving formatting/number formats for the copied cells.
$ws.Range("A3:R3").Copy($ws.Range("A4:R4"))

# Now update row 3 in place with the new weekly values.
$ws.Range("D3").Value = 44425
$ws.Range("J3").Value = 30
$ws.Range("K3").Value = 13000
$ws.Range("L3").Value = 13000
$ws.Range("M3").Value = 13000
$ws.Range("P3").Value = 1300
